$wb = $excel.ActiveWorkbook

# 1) Update "Last Updated" timestamp on the Metadata sheet.
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 11:56 AM"

# 2) Update the "Stock List" sheet: a new record (CAPTRU-RE1) is inserted at
#    row 2, pushing every existing row (2..76) down by one; the former last
#    row (76) falls off the bottom of the table.
$ws = $wb.Worksheets.Item("Stock List")

$firstRow = 2
$lastRow = 76

# Snapshot the current B/C/D/E/H values for rows 2..75 (row 76's original
# data is dropped, so it doesn't need to be read).
$B = @{}
$C = @{}
$D = @{}
$E = @{}
$H = @{}
for ($r = $firstRow; $r -le ($lastRow - 1); $r++) {
    $B[$r] = $ws.Cells.Item($r, 2).Value()
    $C[$r] = $ws.Cells.Item($r, 3).Value()
    $D[$r] = $ws.Cells.Item($r, 4).Value()
    $E[$r] = $ws.Cells.Item($r, 5).Value()
    $H[$r] = $ws.Cells.Item($r, 8).Value()
}

# Shift rows 3..76 to hold what used to be in rows 2..75.
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 2).Value = $B[$src]
    $ws.Cells.Item($r, 3).Value = $C[$src]
    $ws.Cells.Item($r, 4).Value = $D[$src]
    $ws.Cells.Item($r, 5).Value = $E[$src]
    $ws.Cells.Item($r, 8).Value = $H[$src]
}

# Write the brand-new row 2 record.
$ws.Cells.Item($firstRow, 2).Value = "CAPTRU-RE1"
$ws.Cells.Item($firstRow, 3).Value = "CAPTRU-RE1"
$ws.Cells.Item($firstRow, 4).Value = 5.67
$ws.Cells.Item($firstRow, 5).Value = -11.9565
